$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count

for ($r = 1; $r -le $rows; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -ne $null -and $val -like "*POS=Foc*") {
        $cell.Value = $val -replace "POS=Foc", "Other=Foc"
    }
}

$a450 = $ws.Cells.Item(450, 1)
$a450val = $a450.Value2
if ($a450val.StartsWith("*")) {
    $a450.Value = $a450val.Substring(1)
}
